{"js": "// The document contains two true/false question blocks, each with four\n// sub-statements a) b) c) d) that end with an inline answer marker\n// \" (\u0111)\" (true) or \" (s)\" (false). This edit removes those inline\n// markers from the statement text and instead appends a new summary\n// paragraph \"\u0110\u00e1p \u00e1n: XXXX\" (X = \u0110 or S) right after item d) of each\n// block, combining the four answers in order.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Regex that matches the trailing inline answer marker, e.g.\n// \" (\u0111)\", \" (s)\", or the marker glued directly to the sentence \".(s)\".\nconst markerRe = /\\s*\\((\u0111|\u0110|s|S)\\)\\s*$/u;\n\n// Each block: the paragraph holding item d) together with the combined\n// answer string to append right after it.\nconst blocks = [\n  { dText: \"d) T\u1ea1o phim b\u1eb1ng ph\u1ea7n m\u1ec1m y\u00eau c\u1ea7u s\u1eafp x\u1ebfp t\u01b0 li\u1ec7u theo th\u1ee9 t\u1ef1 ng\u1eabu nhi\u00ean \u0111\u1ec3 t\u1ea1o n\u00ean chu\u1ed7i c\u00e1c ph\u00e2n c\u1ea3nh. (s)\", answer: \"\u0110S\u0110S\" },\n  { dText: \"d) Ng\u0103n t\u01b0 li\u1ec7u l\u00e0 n\u01a1i xem tr\u01b0\u1edbc c\u00e1c ph\u00e2n c\u1ea3nh \u0111\u00e3 ho\u00e0n ch\u1ec9nh sau khi t\u1ea1o phim.(s)\", answer: \"\u0110\u0110\u0110S\" },\n];\n\n// Collect the paragraphs that need their trailing marker stripped.\nconst items = paragraphs.items;\nconst toClean = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (markerRe.test(t)) {\n    toClean.push(i);\n  }\n}\n\n// Strip the trailing marker from every matching statement paragraph.\nfor (const idx of toClean) {\n  const p = items[idx];\n  const cleaned = p.text.replace(markerRe, \"\");\n  p.getRange().insertText(cleaned, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Now find the (still valid) paragraphs for each block's item d) by their\n// now-cleaned text, and insert the \"\u0110\u00e1p \u00e1n: ...\" paragraph right after.\nfor (const block of blocks) {\n  const cleanedDText = block.dText.replace(markerRe, \"\");\n  const results = body.search(cleanedDText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate paragraph for: \" + cleanedDText);\n  }\n  const range = results.items[0];\n  const paragraph = range.paragraphs.getFirst();\n  paragraph.insertParagraph(\"\u0110\u00e1p \u00e1n: \" + block.answer, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# The document contains two true/false question blocks, each with four\n# sub-statements a) b) c) d) that end with an inline answer marker\n# \" (\u0111)\" (true) or \" (s)\" (false). This edit removes those inline\n# markers from the statement text and instead appends a new summary\n# paragraph \"\u0110\u00e1p \u00e1n: XXXX\" (X = \u0110 or S) right after item d) of each\n# block, combining the four answers in order.\n\n$d = $word.ActiveDocument\n\n# wdReplaceNone = 0, wdReplaceOne = 1\n$wdReplaceNone = 0\n$wdReplaceOne = 1\n\n# Literal marker substrings to strip (verbatim, including any leading\n# space that is actually present in the source text). One entry per\n# occurrence, in document order. The very last one (\"Ng\u0103n t\u01b0 li\u1ec7u...\")\n# has no space before the marker in the source text.\n$markers = @(\n    \" (\u0111)\",\n    \" (s)\",\n    \" (\u0111)\",\n    \" (s)\",\n    \" (\u0111)\",\n    \" (\u0111)\",\n    \" (\u0111)\",\n    \"(s)\"\n)\n\nforeach ($marker in $markers) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, \"\", $wdReplaceOne) | Out-Null\n}\n\nfunction Insert-AnswerParagraphAfter($searchText, $answerText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, \"\", $wdReplaceNone) | Out-Null\n\n    # Compute the 1-based paragraph index of the found range by counting\n    # how many paragraphs precede its start.\n    $countRange = $d.Range(0, $rng.Start)\n    $idx = $countRange.Paragraphs.Count + 1\n\n    $d.Paragraphs.Item($idx).Range.InsertParagraphAfter()\n    $d.Paragraphs.Item($idx + 1).Range.Text = $answerText\n}\n\n# Block 1 (C\u00e2u 1): item d), now cleaned of its trailing marker.\n$dText1 = \"d) T\u1ea1o phim b\u1eb1ng ph\u1ea7n m\u1ec1m y\u00eau c\u1ea7u s\u1eafp x\u1ebfp t\u01b0 li\u1ec7u theo th\u1ee9 t\u1ef1 ng\u1eabu nhi\u00ean \u0111\u1ec3 t\u1ea1o n\u00ean chu\u1ed7i c\u00e1c ph\u00e2n c\u1ea3nh.\"\nInsert-AnswerParagraphAfter $dText1 \"\u0110\u00e1p \u00e1n: \u0110S\u0110S\"\n\n# Block 2 (C\u00e2u 2): item d), now cleaned of its trailing marker.\n$dText2 = \"d) Ng\u0103n t\u01b0 li\u1ec7u l\u00e0 n\u01a1i xem tr\u01b0\u1edbc c\u00e1c ph\u00e2n c\u1ea3nh \u0111\u00e3 ho\u00e0n ch\u1ec9nh sau khi t\u1ea1o phim.\"\nInsert-AnswerParagraphAfter $dText2 \"\u0110\u00e1p \u00e1n: \u0110\u0110\u0110S\"\n"}
